$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "75"
$ws.Range("C3").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "248120.00"
$ws.Range("D3").Style = "Normal"

$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = "353"
$ws.Range("C6").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "937010.82"
$ws.Range("D6").Style = "Normal"

$ws.Range("C7").NumberFormat = "@"
$ws.Range("C7").Value = "55"
$ws.Range("C7").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "136000.00"
$ws.Range("D7").Style = "Normal"

$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = "697"
$ws.Range("C8").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "2448701.92"
$ws.Range("D8").Style = "Normal"

$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = "23"
$ws.Range("C9").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "57100.00"
$ws.Range("D9").Style = "Normal"

$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = "13"
$ws.Range("C10").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35500.00"
$ws.Range("D10").Style = "Normal"

$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = "33"
$ws.Range("C11").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "77270.00"
$ws.Range("D11").Style = "Normal"

$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = "135"
$ws.Range("C12").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "377142.00"
$ws.Range("D12").Style = "Normal"

$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "84"
$ws.Range("C14").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "209988.98"
$ws.Range("D14").Style = "Normal"

$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "109"
$ws.Range("C16").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "464258.76"
$ws.Range("D16").Style = "Normal"

$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = "166"
$ws.Range("C17").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "380944.87"
$ws.Range("D17").Style = "Normal"

$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").Value = "360"
$ws.Range("C37").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1406311.70"
$ws.Range("D37").Style = "Normal"

$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = "18"
$ws.Range("C40").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "53170.00"
$ws.Range("D40").Style = "Normal"

$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = "75"
$ws.Range("C46").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "169893.00"
$ws.Range("D46").Style = "Normal"

$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = "16"
$ws.Range("C47").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "52500.00"
$ws.Range("D47").Style = "Normal"

$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "112"
$ws.Range("C50").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "286968.33"
$ws.Range("D50").Style = "Normal"

$ws.Range("C52").NumberFormat = "@"
$ws.Range("C52").Value = "265"
$ws.Range("C52").Style = "Normal"
$ws.Range("D52").NumberFormat = "@"
$ws.Range("D52").Value = "969715.67"
$ws.Range("D52").Style = "Normal"

$ws.Range("C60").NumberFormat = "@"
$ws.Range("C60").Value = "42"
$ws.Range("C60").Style = "Normal"
$ws.Range("D60").NumberFormat = "@"
$ws.Range("D60").Value = "174656.00"
$ws.Range("D60").Style = "Normal"

$ws.Range("C74").NumberFormat = "@"
$ws.Range("C74").Value = "17"
$ws.Range("C74").Style = "Normal"
$ws.Range("D74").NumberFormat = "@"
$ws.Range("D74").Value = "65500.00"
$ws.Range("D74").Style = "Normal"

$ws.Range("C75").NumberFormat = "@"
$ws.Range("C75").Value = "40"
$ws.Range("C75").Style = "Normal"
$ws.Range("D75").NumberFormat = "@"
$ws.Range("D75").Value = "114579.25"
$ws.Range("D75").Style = "Normal"

$ws.Range("C77").NumberFormat = "@"
$ws.Range("C77").Value = "86"
$ws.Range("C77").Style = "Normal"
$ws.Range("D77").NumberFormat = "@"
$ws.Range("D77").Value = "227487.00"
$ws.Range("D77").Style = "Normal"

$ws.Range("C78").NumberFormat = "@"
$ws.Range("C78").Value = "199"
$ws.Range("C78").Style = "Normal"
$ws.Range("D78").NumberFormat = "@"
$ws.Range("D78").Value = "543693.00"
$ws.Range("D78").Style = "Normal"

$ws.Range("C84").NumberFormat = "@"
$ws.Range("C84").Value = "69"
$ws.Range("C84").Style = "Normal"
$ws.Range("D84").NumberFormat = "@"
$ws.Range("D84").Value = "228657.55"
$ws.Range("D84").Style = "Normal"

$ws.Range("C85").NumberFormat = "@"
$ws.Range("C85").Value = "33"
$ws.Range("C85").Style = "Normal"
$ws.Range("D85").NumberFormat = "@"
$ws.Range("D85").Value = "104669.00"
$ws.Range("D85").Style = "Normal"

$ws.Range("C86").NumberFormat = "@"
$ws.Range("C86").Value = "43"
$ws.Range("C86").Style = "Normal"
$ws.Range("D86").NumberFormat = "@"
$ws.Range("D86").Value = "99500.00"
$ws.Range("D86").Style = "Normal"

$ws.Range("C106").NumberFormat = "@"
$ws.Range("C106").Value = "21"
$ws.Range("C106").Style = "Normal"
$ws.Range("D106").NumberFormat = "@"
$ws.Range("D106").Value = "58209.84"
$ws.Range("D106").Style = "Normal"

$ws.Range("C107").NumberFormat = "@"
$ws.Range("C107").Value = "70"
$ws.Range("C107").Style = "Normal"
$ws.Range("D107").NumberFormat = "@"
$ws.Range("D107").Value = "174310.00"
$ws.Range("D107").Style = "Normal"

$ws.Range("C108").NumberFormat = "@"
$ws.Range("C108").Value = "34"
$ws.Range("C108").Style = "Normal"
$ws.Range("D108").NumberFormat = "@"
$ws.Range("D108").Value = "105434.00"
$ws.Range("D108").Style = "Normal"

$ws.Range("C110").NumberFormat = "@"
$ws.Range("C110").Value = "83"
$ws.Range("C110").Style = "Normal"
$ws.Range("D110").NumberFormat = "@"
$ws.Range("D110").Value = "525071.82"
$ws.Range("D110").Style = "Normal"

$ws.Range("C112").NumberFormat = "@"
$ws.Range("C112").Value = "4"
$ws.Range("C112").Style = "Normal"
$ws.Range("D112").NumberFormat = "@"
$ws.Range("D112").Value = "9500.00"
$ws.Range("D112").Style = "Normal"

$ws.Range("C113").NumberFormat = "@"
$ws.Range("C113").Value = "23"
$ws.Range("C113").Style = "Normal"
$ws.Range("D113").NumberFormat = "@"
$ws.Range("D113").Value = "70767.00"
$ws.Range("D113").Style = "Normal"

$ws.Range("C114").NumberFormat = "@"
$ws.Range("C114").Value = "25"
$ws.Range("C114").Style = "Normal"
$ws.Range("D114").NumberFormat = "@"
$ws.Range("D114").Value = "69895.00"
$ws.Range("D114").Style = "Normal"

$ws.Range("C121").NumberFormat = "@"
$ws.Range("C121").Value = "63"
$ws.Range("C121").Style = "Normal"
$ws.Range("D121").NumberFormat = "@"
$ws.Range("D121").Value = "172877.00"
$ws.Range("D121").Style = "Normal"

$ws.Range("C122").NumberFormat = "@"
$ws.Range("C122").Value = "240"
$ws.Range("C122").Style = "Normal"
$ws.Range("D122").NumberFormat = "@"
$ws.Range("D122").Value = "655508.00"
$ws.Range("D122").Style = "Normal"

$ws.Range("C123").NumberFormat = "@"
$ws.Range("C123").Value = "103"
$ws.Range("C123").Style = "Normal"
$ws.Range("D123").NumberFormat = "@"
$ws.Range("D123").Value = "274081.45"
$ws.Range("D123").Style = "Normal"

$ws.Range("C124").NumberFormat = "@"
$ws.Range("C124").Value = "467"
$ws.Range("C124").Style = "Normal"
$ws.Range("D124").NumberFormat = "@"
$ws.Range("D124").Value = "2053132.06"
$ws.Range("D124").Style = "Normal"

$ws.Range("C128").NumberFormat = "@"
$ws.Range("C128").Value = "87"
$ws.Range("C128").Style = "Normal"
$ws.Range("D128").NumberFormat = "@"
$ws.Range("D128").Value = "262743.68"
$ws.Range("D128").Style = "Normal"

$ws.Range("C132").NumberFormat = "@"
$ws.Range("C132").Value = "83"
$ws.Range("C132").Style = "Normal"
$ws.Range("D132").NumberFormat = "@"
$ws.Range("D132").Value = "369163.75"
$ws.Range("D132").Style = "Normal"

$ws.Range("C135").NumberFormat = "@"
$ws.Range("C135").Value = "207"
$ws.Range("C135").Style = "Normal"
$ws.Range("D135").NumberFormat = "@"
$ws.Range("D135").Value = "576620.00"
$ws.Range("D135").Style = "Normal"

$ws.Range("C138").NumberFormat = "@"
$ws.Range("C138").Value = "562"
$ws.Range("C138").Style = "Normal"
$ws.Range("D138").NumberFormat = "@"
$ws.Range("D138").Value = "1406546.00"
$ws.Range("D138").Style = "Normal"

$ws.Range("C139").NumberFormat = "@"
$ws.Range("C139").Value = "1793"
$ws.Range("C139").Style = "Normal"
$ws.Range("D139").NumberFormat = "@"
$ws.Range("D139").Value = "4814610.93"
$ws.Range("D139").Style = "Normal"

$ws.Range("C140").NumberFormat = "@"
$ws.Range("C140").Value = "2480"
$ws.Range("C140").Style = "Normal"
$ws.Range("D140").NumberFormat = "@"
$ws.Range("D140").Value = "6201465.94"
$ws.Range("D140").Style = "Normal"

$ws.Range("C141").NumberFormat = "@"
$ws.Range("C141").Value = "2499"
$ws.Range("C141").Style = "Normal"
$ws.Range("D141").NumberFormat = "@"
$ws.Range("D141").Value = "10540285.14"
$ws.Range("D141").Style = "Normal"

$ws.Range("C142").NumberFormat = "@"
$ws.Range("C142").Value = "352"
$ws.Range("C142").Style = "Normal"
$ws.Range("D142").NumberFormat = "@"
$ws.Range("D142").Value = "990954.51"
$ws.Range("D142").Style = "Normal"

$ws.Range("C143").NumberFormat = "@"
$ws.Range("C143").Value = "124"
$ws.Range("C143").Style = "Normal"
$ws.Range("D143").NumberFormat = "@"
$ws.Range("D143").Value = "304500.00"
$ws.Range("D143").Style = "Normal"

$ws.Range("C144").NumberFormat = "@"
$ws.Range("C144").Value = "240"
$ws.Range("C144").Style = "Normal"
$ws.Range("D144").NumberFormat = "@"
$ws.Range("D144").Value = "610800.00"
$ws.Range("D144").Style = "Normal"

$ws.Range("C145").NumberFormat = "@"
$ws.Range("C145").Value = "1017"
$ws.Range("C145").Style = "Normal"
$ws.Range("D145").NumberFormat = "@"
$ws.Range("D145").Value = "2642933.25"
$ws.Range("D145").Style = "Normal"

$ws.Range("C146").NumberFormat = "@"
$ws.Range("C146").Value = "486"
$ws.Range("C146").Style = "Normal"
$ws.Range("D146").NumberFormat = "@"
$ws.Range("D146").Value = "1402044.49"
$ws.Range("D146").Style = "Normal"

$ws.Range("C147").NumberFormat = "@"
$ws.Range("C147").Value = "368"
$ws.Range("C147").Style = "Normal"
$ws.Range("D147").NumberFormat = "@"
$ws.Range("D147").Value = "923700.16"
$ws.Range("D147").Style = "Normal"

$ws.Range("C148").NumberFormat = "@"
$ws.Range("C148").Value = "147"
$ws.Range("C148").Style = "Normal"
$ws.Range("D148").NumberFormat = "@"
$ws.Range("D148").Value = "361000.00"
$ws.Range("D148").Style = "Normal"

$ws.Range("C149").NumberFormat = "@"
$ws.Range("C149").Value = "402"
$ws.Range("C149").Style = "Normal"
$ws.Range("D149").NumberFormat = "@"
$ws.Range("D149").Value = "1252348.46"
$ws.Range("D149").Style = "Normal"

$ws.Range("C150").NumberFormat = "@"
$ws.Range("C150").Value = "841"
$ws.Range("C150").Style = "Normal"
$ws.Range("D150").NumberFormat = "@"
$ws.Range("D150").Value = "2026630.82"
$ws.Range("D150").Style = "Normal"

$ws.Range("C194").NumberFormat = "@"
$ws.Range("C194").Value = "54"
$ws.Range("C194").Style = "Normal"
$ws.Range("D194").NumberFormat = "@"
$ws.Range("D194").Value = "166300.00"
$ws.Range("D194").Style = "Normal"

$ws.Range("C195").NumberFormat = "@"
$ws.Range("C195").Value = "55"
$ws.Range("C195").Style = "Normal"
$ws.Range("D195").NumberFormat = "@"
$ws.Range("D195").Value = "139153.50"
$ws.Range("D195").Style = "Normal"

$ws.Range("C197").NumberFormat = "@"
$ws.Range("C197").Value = "350"
$ws.Range("C197").Style = "Normal"
$ws.Range("D197").NumberFormat = "@"
$ws.Range("D197").Value = "945788.00"
$ws.Range("D197").Style = "Normal"

$ws.Range("C199").NumberFormat = "@"
$ws.Range("C199").Value = "644"
$ws.Range("C199").Style = "Normal"
$ws.Range("D199").NumberFormat = "@"
$ws.Range("D199").Value = "2404258.16"
$ws.Range("D199").Style = "Normal"

$ws.Range("C200").NumberFormat = "@"
$ws.Range("C200").Value = "25"
$ws.Range("C200").Style = "Normal"
$ws.Range("D200").NumberFormat = "@"
$ws.Range("D200").Value = "96238.00"
$ws.Range("D200").Style = "Normal"

$ws.Range("C202").NumberFormat = "@"
$ws.Range("C202").Value = "32"
$ws.Range("C202").Style = "Normal"
$ws.Range("D202").NumberFormat = "@"
$ws.Range("D202").Value = "82000.00"
$ws.Range("D202").Style = "Normal"

$ws.Range("C203").NumberFormat = "@"
$ws.Range("C203").Value = "152"
$ws.Range("C203").Style = "Normal"
$ws.Range("D203").NumberFormat = "@"
$ws.Range("D203").Value = "449133.00"
$ws.Range("D203").Style = "Normal"

$ws.Range("C205").NumberFormat = "@"
$ws.Range("C205").Value = "78"
$ws.Range("C205").Style = "Normal"
$ws.Range("D205").NumberFormat = "@"
$ws.Range("D205").Value = "183005.00"
$ws.Range("D205").Style = "Normal"

$ws.Range("C207").NumberFormat = "@"
$ws.Range("C207").Value = "119"
$ws.Range("C207").Style = "Normal"
$ws.Range("D207").NumberFormat = "@"
$ws.Range("D207").Value = "553095.50"
$ws.Range("D207").Style = "Normal"

$ws.Range("C208").NumberFormat = "@"
$ws.Range("C208").Value = "137"
$ws.Range("C208").Style = "Normal"
$ws.Range("D208").NumberFormat = "@"
$ws.Range("D208").Value = "304196.77"
$ws.Range("D208").Style = "Normal"

